# Apply targeted cell additions to Sheet1, matching the target diff.
# The diff only adds brand-new values into previously-empty cells; no
# existing cell's displayed value changes (the shared-string index churn
# in the XML diff is just an artifact of new strings being inserted
# into the shared string table).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# F1 is an untouched blank cell in the source file; explicitly keep it blank
# so that the COM round-trip doesn't let it pick up a stray value.
$ws.Range("F1").Value = ""

# Practice rows (2-5): fill in column D with the carrier word for each practice trial.
$ws.Range("D2").Value = "can"
$ws.Range("D3").Value = "where"
$ws.Range("D4").Value = "do"
$ws.Range("D5").Value = "look"

# Generic rows (6-9): fill in column J with pair_kind for the new unique_video/unique_audio pairs.
$ws.Range("J6").Value = "unique_video"
$ws.Range("J7").Value = "unique_video"
$ws.Range("J8").Value = "unique_audio"
$ws.Range("J9").Value = "unique_audio"

# New unique_video / unique_audio rows (14-21): fill in kind (C) and carrier (D).
$ws.Range("C14").Value = "unique_video"
$ws.Range("D14").Value = "can"

$ws.Range("C15").Value = "unique_video"
$ws.Range("D15").Value = "can"

$ws.Range("C16").Value = "unique_video"
$ws.Range("D16").Value = "do"

$ws.Range("C17").Value = "unique_video"
$ws.Range("D17").Value = "do"

$ws.Range("C18").Value = "unique_audio"
$ws.Range("D18").Value = "look"

$ws.Range("C19").Value = "unique_audio"
$ws.Range("D19").Value = "look"

$ws.Range("C20").Value = "unique_audio"
$ws.Range("D20").Value = "where"

$ws.Range("C21").Value = "unique_audio"
$ws.Range("D21").Value = "where"
